# Automatische test-sync: 2025-06-19 19:06:30
# Appends a new incoming-mail log row to the "Logs" sheet and refreshes the
# "Dashboard" category-count sheet to match.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 37 -----------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A37").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B37").Value = "mailmind.test@zohomail.eu"
$logs.Range("C37").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D37").Value = "Offerte-aanvraag"
$logs.Range("F37").Value = "2025-06-19 19:06:25"
$logs.Range("G37").Value = "Nee"

# Extend the two conditional-formatting blocks (Categorie + Beantwoord
# columns) so they keep covering the whole data range through row 37.
$dCond = $logs.Range("D2:D36").FormatConditions.Item(1)
$dCond.ModifyAppliesToRange($logs.Range("D2:D37"))

$gCond = $logs.Range("G2:G36").FormatConditions.Item(1)
$gCond.ModifyAppliesToRange($logs.Range("G2:G37"))

# --- Dashboard sheet: update the Offerte-aanvraag / Openingstijden rows
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Offerte-aanvraag"
$dash.Range("B7").Value = 3
$dash.Range("A8").Value = "Openingstijden"
$dash.Range("B8").Value = 2
